$d = $word.ActiveDocument

# 1) Title paragraph: merge the word-by-word runs into a single run of text.
$d.Content.Find.Execute(
    "Answers: Introduction to differentiation and the derivative",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answers: Introduction to differentiation and the derivative", 2)

# 2) Author paragraph: merge the word-by-word runs into a single run of text.
$d.Content.Find.Execute(
    "Sara Delgado Garcia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sara Delgado Garcia", 2)

# 3) Abstract paragraph: merge the word-by-word runs into a single run of text.
$d.Content.Find.Execute(
    "Answers to questions relating to the guide on introduction to differentiation and the derivative.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Answers to questions relating to the guide on introduction to differentiation and the derivative.", 2)
